$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.774.04"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.359.79"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.10"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "660.34"
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("E7").Value = "  -2.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.428"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -3.91%  "
$ws.Range("D11").Value = "3.357.16"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.04"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "97.518.36"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.13"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").Value = "3.987.31"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("E18").Value = "  +3.91%  "
$ws.Range("D19").Value = "3.381.63"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("E20").Value = "  +3.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.530"
$ws.Range("E21").Value = "  +2.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.89"
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "513.47"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.96"
$ws.Range("E26").Value = "  +14.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.85"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.42"
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.147"
$ws.Range("E29").Value = "  -5.04%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.63"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("B32").Value = "Cronos"
$ws.Range("C32").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.191"
$ws.Range("E32").Value = "  -7.22%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.62"
$ws.Range("E33").Value = "  +15.51%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("B35").Value = "PolygonEcosystemToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.570"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.78"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.52"
$ws.Range("E37").Value = "  +7.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.97"
$ws.Range("E38").Value = "  +4.38%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "525.11"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.153"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0442"
$ws.Range("E42").Value = "  +5.07%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.42"
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.862"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.76"
$ws.Range("E45").Value = "  +10.35%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.69"
$ws.Range("E46").Value = "  +5.94%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.73"
$ws.Range("E47").Value = "  +6.27%  "
$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.63"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.79"
$ws.Range("E49").Value = "  +5.17%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.16"
$ws.Range("E50").Value = "  -3.21%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.06"
$ws.Range("E51").Value = "  -0.12%  "
